$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.593.75'
$ws.Range("E2").Value = '  +3.09%  '

$ws.Range("D3").Value = '4.006.10'
$ws.Range("E3").Value = '  +1.55%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '531.43'
$ws.Range("E5").Value = '  +5.64%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.10'
$ws.Range("E6").Value = '  +0.78%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.624'
$ws.Range("E7").Value = '  -0.33%  '

$ws.Range("E8").Value = '  +0.14%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.737'
$ws.Range("E9").Value = '  +0.32%  '

$ws.Range("E10").Value = '  +0.57%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000347'
$ws.Range("E11").Value = '  -1.18%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.15'
$ws.Range("E12").Value = '  -1.30%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.64'
$ws.Range("E13").Value = '  +1.19%  '

$ws.Range("D14").Value = '4.640.81'
$ws.Range("E14").Value = '  +1.47%  '

$ws.Range("D15").Value = '4.011.44'
$ws.Range("E15").Value = '  +1.83%  '

$ws.Range("B16").Value = 'Uniswap'
$ws.Range("C16").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.33'
$ws.Range("E16").Value = '  +0.47%  '

$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '21.17'
$ws.Range("E17").Value = '  +5.61%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.23'
$ws.Range("E18").Value = '  +2.51%  '

$ws.Range("E19").Value = '  -1.95%  '

$ws.Range("D20").Value = '71.411.52'
$ws.Range("E20").Value = '  +2.87%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '444.23'
$ws.Range("E21").Value = '  +1.58%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.54'
$ws.Range("E22").Value = '  +2.44%  '

$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '91.87'
$ws.Range("E23").Value = '  +3.23%  '

$ws.Range("B24").Value = 'RenderToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.29'
$ws.Range("E24").Value = '  +1.92%  '

$ws.Range("E25").Value = '  -3.62%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.10'
$ws.Range("E26").Value = '  +5.53%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.78'
$ws.Range("E27").Value = '  -3.88%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.05'

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '13.58'
$ws.Range("E29").Value = '  +0.71%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '694.22'
$ws.Range("E30").Value = '  -2.69%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.129'
$ws.Range("E31").Value = '  -0.35%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.89'
$ws.Range("E32").Value = '  -0.15%  '

$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '68.87'
$ws.Range("E33").Value = '  +7.13%  '

$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.74'
$ws.Range("E34").Value = '  +11.16%  '

$ws.Range("D35").Value = '0.0₃0910'
$ws.Range("E35").Value = '  +1.98%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.441'
$ws.Range("E36").Value = '  -2.92%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '40.81'
$ws.Range("E37").Value = '  -0.69%  '

$ws.Range("B38").Value = 'ThetaToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.54'
$ws.Range("E38").Value = '  +14.53%  '

$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.151'
$ws.Range("E39").Value = '  -0.39%  '

$ws.Range("E40").Value = '  -0.04%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.22%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0491'
$ws.Range("E42").Value = '  -0.17%  '

$ws.Range("E43").Value = '  -0.15%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.10'
$ws.Range("E44").Value = '  +1.26%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.37'
$ws.Range("E45").Value = '  +11.43%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.48'
$ws.Range("E46").Value = '  +2.87%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.144'
$ws.Range("E47").Value = '  -0.16%  '

$ws.Range("D48").Value = '0.0₆0364'
$ws.Range("E48").Value = '  +4.49%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000284'
$ws.Range("E49").Value = '  +18.76%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.32'
$ws.Range("E50").Value = '  +5.48%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.40'
$ws.Range("E51").Value = '  -0.06%  '
